# Natmi LR-pairs (Myoc-Fzd4) results refreshed following Dr Hou's advice:
# the cluster set grew from {FAPs, sCs} to {ECs, FAPs, sCs}, so the
# Sending/Target cluster cross-product goes from 2x2=4 to 3x3=9 result rows,
# and every metric column is recomputed against the enlarged data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Myoc/Fzd4 -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Myoc"
$ws.Range("C2").Value = "Fzd4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.349434
$ws.Range("H2").Value = 1.048302
$ws.Range("I2").Value = 0.0009963999680650763
$ws.Range("J2").Value = 0.0009963999680650763
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 22.41709
$ws.Range("N2").Value = 67.25127000000001
$ws.Range("O2").Value = 0.3988455747018376
$ws.Range("P2").Value = 0.3988455747018376
$ws.Range("Q2").Value = 7.833293427060001
$ws.Range("R2").Value = 70.49964084354001
$ws.Range("S2").Value = 0.0003974097178958079
$ws.Range("T2").Value = 0.000397409717895808

# Row 3: ECs -> Myoc/Fzd4 -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Myoc"
$ws.Range("C3").Value = "Fzd4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.349434
$ws.Range("H3").Value = 1.048302
$ws.Range("I3").Value = 0.0009963999680650763
$ws.Range("J3").Value = 0.0009963999680650763
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 16.78189033333333
$ws.Range("N3").Value = 50.345671
$ws.Range("O3").Value = 0.2985839238983091
$ws.Range("P3").Value = 0.2985839238983091
$ws.Range("Q3").Value = 5.864163066738
$ws.Range("R3").Value = 52.777467600642
$ws.Range("S3").Value = 0.0002975090122370203
$ws.Range("T3").Value = 0.0002975090122370204

# Row 4: ECs -> Myoc/Fzd4 -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Myoc"
$ws.Range("C4").Value = "Fzd4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.349434
$ws.Range("H4").Value = 1.048302
$ws.Range("I4").Value = 0.0009963999680650763
$ws.Range("J4").Value = 0.0009963999680650763
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 17.00595566666667
$ws.Range("N4").Value = 51.017867
$ws.Range("O4").Value = 0.3025705013998533
$ws.Range("P4").Value = 0.3025705013998533
$ws.Range("Q4").Value = 5.942459112426
$ws.Range("R4").Value = 53.482132011834
$ws.Range("S4").Value = 0.0003014812379322479
$ws.Range("T4").Value = 0.0003014812379322479

# Row 5: FAPs -> Myoc/Fzd4 -> ECs  (new row)
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Myoc"
$ws.Range("C5").Value = "Fzd4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 348.977468
$ws.Range("H5").Value = 1046.932404
$ws.Range("I5").Value = 0.9950981815468188
$ws.Range("J5").Value = 0.9950981815468188
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 22.41709
$ws.Range("N5").Value = 67.25127000000001
$ws.Range("O5").Value = 0.3988455747018376
$ws.Range("P5").Value = 0.3988455747018376
$ws.Range("Q5").Value = 7823.05930812812
$ws.Range("R5").Value = 70407.53377315307
$ws.Range("S5").Value = 0.3968905061037944
$ws.Range("T5").Value = 0.3968905061037945

# Row 6: FAPs -> Myoc/Fzd4 -> FAPs  (new row)
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Myoc"
$ws.Range("C6").Value = "Fzd4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 348.977468
$ws.Range("H6").Value = 1046.932404
$ws.Range("I6").Value = 0.9950981815468188
$ws.Range("J6").Value = 0.9950981815468188
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 16.78189033333333
$ws.Range("N6").Value = 50.345671
$ws.Range("O6").Value = 0.2985839238983091
$ws.Range("P6").Value = 0.2985839238983091
$ws.Range("Q6").Value = 5856.501596780342
$ws.Range("R6").Value = 52708.51437102308
$ws.Range("S6").Value = 0.2971203197103211
$ws.Range("T6").Value = 0.2971203197103212

# Row 7: FAPs -> Myoc/Fzd4 -> sCs  (new row)
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Myoc"
$ws.Range("C7").Value = "Fzd4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 348.977468
$ws.Range("H7").Value = 1046.932404
$ws.Range("I7").Value = 0.9950981815468188
$ws.Range("J7").Value = 0.9950981815468188
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 17.00595566666667
$ws.Range("N7").Value = 51.017867
$ws.Range("O7").Value = 0.3025705013998533
$ws.Range("P7").Value = 0.3025705013998533
$ws.Range("Q7").Value = 5934.695349473584
$ws.Range("R7").Value = 53412.25814526226
$ws.Range("S7").Value = 0.3010873557327032
$ws.Range("T7").Value = 0.3010873557327032

# Row 8: sCs -> Myoc/Fzd4 -> ECs  (new row)
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Myoc"
$ws.Range("C8").Value = "Fzd4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.369616666666667
$ws.Range("H8").Value = 4.10885
$ws.Range("I8").Value = 0.00390541848511611
$ws.Range("J8").Value = 0.00390541848511611
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 22.41709
$ws.Range("N8").Value = 67.25127000000001
$ws.Range("O8").Value = 0.3988455747018376
$ws.Range("P8").Value = 0.3988455747018376
$ws.Range("Q8").Value = 30.70282008216667
$ws.Range("R8").Value = 276.3253807395
$ws.Range("S8").Value = 0.001557658880147315
$ws.Range("T8").Value = 0.001557658880147315

# Row 9: sCs -> Myoc/Fzd4 -> FAPs  (new row)
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Myoc"
$ws.Range("C9").Value = "Fzd4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.369616666666667
$ws.Range("H9").Value = 4.10885
$ws.Range("I9").Value = 0.00390541848511611
$ws.Range("J9").Value = 0.00390541848511611
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 16.78189033333333
$ws.Range("N9").Value = 50.345671
$ws.Range("O9").Value = 0.2985839238983091
$ws.Range("P9").Value = 0.2985839238983091
$ws.Range("Q9").Value = 22.98475669870556
$ws.Range("R9").Value = 206.86281028835
$ws.Range("S9").Value = 0.001166095175750958
$ws.Range("T9").Value = 0.001166095175750958

# Row 10: sCs -> Myoc/Fzd4 -> sCs  (new row)
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Myoc"
$ws.Range("C10").Value = "Fzd4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.369616666666667
$ws.Range("H10").Value = 4.10885
$ws.Range("I10").Value = 0.00390541848511611
$ws.Range("J10").Value = 0.00390541848511611
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 17.00595566666667
$ws.Range("N10").Value = 51.017867
$ws.Range("O10").Value = 0.3025705013998533
$ws.Range("P10").Value = 0.3025705013998533
$ws.Range("Q10").Value = 23.29164031366111
$ws.Range("R10").Value = 209.62476282295
$ws.Range("S10").Value = 0.001181664429217837
$ws.Range("T10").Value = 0.001181664429217837
